# Fix mismatched PERSON placeholder numbers inside the inflected-form examples
# („o ...“, „s ...“ pairs) so that all three occurrences in each bullet refer to
# the same person, and renumber following the correction.
$d = $word.ActiveDocument

$pairs = @(
    @("[[PERSON_85]] – „o [[PERSON_86]]“, „s [[PERSON_86]]“", "[[PERSON_85]] – „o [[PERSON_85]]“, „s [[PERSON_85]]“"),
    @("[[PERSON_87]] – „o [[PERSON_87]]“, „s [[PERSON_87]]“", "[[PERSON_86]] – „o [[PERSON_86]]“, „s [[PERSON_86]]“"),
    @("[[PERSON_88]] – „s [[PERSON_88]]“, „o [[PERSON_88]]“", "[[PERSON_87]] – „s [[PERSON_87]]“, „o [[PERSON_87]]“"),
    @("[[PERSON_89]] – „o [[PERSON_89]]“, „s [[PERSON_89]]“", "[[PERSON_88]] – „o [[PERSON_88]]“, „s [[PERSON_88]]“"),
    @("[[PERSON_90]] – „s [[PERSON_90]]“, „o [[PERSON_90]]“", "[[PERSON_89]] – „s [[PERSON_89]]“, „o [[PERSON_89]]“"),
    @("[[PERSON_91]] – „o [[PERSON_91]]“, „s [[PERSON_91]]“", "[[PERSON_90]] – „o [[PERSON_90]]“, „s [[PERSON_90]]“"),
    @("[[PERSON_92]] – „s [[PERSON_92]]“, „o [[PERSON_92]]“", "[[PERSON_91]] – „s [[PERSON_91]]“, „o [[PERSON_91]]“"),
    @("[[PERSON_93]] – „o [[PERSON_93]]“, „s [[PERSON_93]]“", "[[PERSON_92]] – „o [[PERSON_92]]“, „s [[PERSON_92]]“"),
    @("[[PERSON_94]] – „s [[PERSON_94]]“, „o [[PERSON_95]]“", "[[PERSON_93]] – „s [[PERSON_94]]“, „o [[PERSON_94]]“"),
    @("[[PERSON_96]] – „o [[PERSON_97]]“, „s [[PERSON_96]]“", "[[PERSON_95]] – „o [[PERSON_96]]“, „s [[PERSON_95]]“"),
    @("[[PERSON_98]] – „s [[PERSON_98]]“, „o [[PERSON_98]]“", "[[PERSON_97]] – „s [[PERSON_97]]“, „o [[PERSON_97]]“"),
    @("[[PERSON_99]] – „o [[PERSON_100]]“, „s [[PERSON_99]]“", "[[PERSON_98]] – „o [[PERSON_99]]“, „s [[PERSON_99]]“"),
    @("[[PERSON_101]] – „s [[PERSON_101]]“, „o [[PERSON_101]]“", "[[PERSON_100]] – „s [[PERSON_100]]“, „o [[PERSON_100]]“"),
    @("[[PERSON_102]] – „s [[PERSON_102]]“, „o [[PERSON_103]]“", "[[PERSON_101]] – „s [[PERSON_101]]“, „o [[PERSON_102]]“"),
    @("[[PERSON_104]] – „s [[PERSON_104]]“, „o [[PERSON_104]]“", "[[PERSON_103]] – „s [[PERSON_103]]“, „o [[PERSON_103]]“"),
    @("[[PERSON_105]] – „o [[PERSON_105]]“, „s [[PERSON_105]]“", "[[PERSON_104]] – „o [[PERSON_104]]“, „s [[PERSON_104]]“"),
    @("[[PERSON_106]] – „s [[PERSON_107]]“, „o [[PERSON_108]]“", "[[PERSON_105]] – „s [[PERSON_106]]“, „o [[PERSON_107]]“"),
    @("[[PERSON_109]] – „s [[PERSON_109]]“, „o [[PERSON_109]]“", "[[PERSON_108]] – „s [[PERSON_108]]“, „o [[PERSON_108]]“"),
    @("[[PERSON_110]] – „o [[PERSON_110]]“, „s [[PERSON_110]]“", "[[PERSON_109]] – „o [[PERSON_109]]“, „s [[PERSON_109]]“"),
    @("[[PERSON_111]] – „s [[PERSON_112]]“, „o [[PERSON_112]]“", "[[PERSON_110]] – „s [[PERSON_111]]“, „o [[PERSON_111]]“"),
    @("[[PERSON_113]] – „o [[PERSON_113]]“, „s [[PERSON_113]]“", "[[PERSON_112]] – „o [[PERSON_112]]“, „s [[PERSON_112]]“"),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done."
